$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new columns before column D ("D:E"), shifting the existing
#    quarterly data (old D..K) to the right (new F..M). This matches the
#    source diff, which adds two new quarters (2018-12-31 and 2018-09-30)
#    ahead of the existing quarters.
# ---------------------------------------------------------------------------
$ws.Columns("D:E").Insert()

# Column F used to be column D before the insert, so it already carries the
# correct number formats (date format for the header rows, #,##0 for the
# data rows). Copy those formats across into the two new columns D:E so the
# new cells look like the rest of the table. Restrict the copy to the sheet's
# used rows (5-102) so we don't balloon the sheet's used range down to the
# bottom of the worksheet.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Populate the two new columns with the new quarterly data.
# ---------------------------------------------------------------------------

# --- Income statement (header row 7, data rows 8-35) ---
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373

$ws.Range("D8").Value2 = 360100
$ws.Range("E8").Value2 = 346100
$ws.Range("D9").Value2 = 88000
$ws.Range("E9").Value2 = 89000
$ws.Range("D10").Value2 = 272100
$ws.Range("E10").Value2 = 257100

$ws.Range("D12").Value2 = 65800
$ws.Range("E12").Value2 = 65500
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = 2900
$ws.Range("E14").Value2 = 4700
$ws.Range("D15").Value2 = 500
$ws.Range("E15").Value2 = 500

$ws.Range("D17").Value2 = 298300
$ws.Range("E17").Value2 = 300100
$ws.Range("D18").Value2 = 61800
$ws.Range("E18").Value2 = 46000

$ws.Range("D20").Value2 = 1600
$ws.Range("E20").Value2 = 2400
$ws.Range("D21").Value2 = 80400
$ws.Range("E21").Value2 = 67000
$ws.Range("D22").Value2 = 0
$ws.Range("E22").Value2 = 0
$ws.Range("D23").Value2 = 63500
$ws.Range("E23").Value2 = 48400
$ws.Range("D24").Value2 = 9300
$ws.Range("E24").Value2 = 7000
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 54100
$ws.Range("E26").Value2 = 41400
$ws.Range("D27").Value2 = 54100
$ws.Range("E27").Value2 = 41400
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0

$ws.Range("D29").Value2 = 2400
$ws.Range("E29").Value2 = 1800

$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = -1600
$ws.Range("E32").Value2 = -2400
$ws.Range("D33").Value2 = 56500
$ws.Range("E33").Value2 = 43200
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 56500
$ws.Range("E35").Value2 = 43200

# --- Balance sheet (header row 38, data rows 41-77) ---
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373

$ws.Range("D41").Value2 = 259400
$ws.Range("E41").Value2 = 311400
$ws.Range("D42").Value2 = 271400
$ws.Range("E42").Value2 = 171000
$ws.Range("D43").Value2 = 243000
$ws.Range("E43").Value2 = 239500
$ws.Range("D44").Value2 = 194100
$ws.Range("E44").Value2 = 192400
$ws.Range("D45").Value2 = 54300
$ws.Range("E45").Value2 = 62400
$ws.Range("D46").Value2 = 1022200
$ws.Range("E46").Value2 = 976700
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("D48").Value2 = 245200
$ws.Range("E48").Value2 = 245900
$ws.Range("D49").Value2 = 375300
$ws.Range("E49").Value2 = 379200
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 28500
$ws.Range("E52").Value2 = 28100
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 1671200
$ws.Range("E54").Value2 = 1629900

$ws.Range("D57").Value2 = 48400
$ws.Range("E57").Value2 = 51300
$ws.Range("D58").Value2 = 0
$ws.Range("E58").Value2 = 0
$ws.Range("D59").Value2 = 234600
$ws.Range("E59").Value2 = 224700
$ws.Range("D60").Value2 = 283000
$ws.Range("E60").Value2 = 276100
$ws.Range("D61").Value2 = 0
$ws.Range("E61").Value2 = 0
$ws.Range("D62").Value2 = 149900
$ws.Range("E62").Value2 = 156900
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 432900
$ws.Range("E66").Value2 = 433000

$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 356400
$ws.Range("E72").Value2 = 329300
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 1238400
$ws.Range("E76").Value2 = 1196900
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0

# --- Cash flow statement (header row 80, data rows 81-102) ---
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373

$ws.Range("D81").Value2 = 56500
$ws.Range("E81").Value2 = 43200

$ws.Range("D83").Value2 = 16900
$ws.Range("E83").Value2 = 18600
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 87200
$ws.Range("E89").Value2 = 88500

$ws.Range("D91").Value2 = -7300
$ws.Range("E91").Value2 = -7600
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -114600
$ws.Range("E94").Value2 = 29300

$ws.Range("D96").Value2 = -30500
$ws.Range("E96").Value2 = -30500
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -24200
$ws.Range("E100").Value2 = -22700
$ws.Range("D101").Value2 = -400
$ws.Range("E101").Value2 = -1300
$ws.Range("D102").Value2 = -52000
$ws.Range("E102").Value2 = 93800

# ---------------------------------------------------------------------------
# 3. A few rows were also restated (not simply shifted right); overwrite the
#    shifted-in cells in those rows with their corrected values.
# ---------------------------------------------------------------------------

# Row 12 - Research Development: columns F,G got new (not shifted) figures.
$ws.Range("F12").Value2 = 66000
$ws.Range("G12").Value2 = 61600

# Row 14 - Non Recurring: columns F,G got new (not shifted) figures.
$ws.Range("F14").Value2 = 4500
$ws.Range("G14").Value2 = 2500

# Row 91 - Capital Expenditures: columns F-J were restated.
$ws.Range("F91").Value2 = -11600
$ws.Range("G91").Value2 = -8100
$ws.Range("H91").Value2 = -6200
$ws.Range("I91").Value2 = -8400
$ws.Range("J91").Value2 = -4900
